# "Plans and Addons xpath update"
# Adds a PlansAndAddonsUsername/Password pair of rows to the LoginPage
# (xpaths) sheet, widens its key column to fit the longer key name, and
# moves the active-sheet/selection state around (SetUp -> B18,
# LoginPage -> B27, MigrationPage becomes the active tab).

$wb = $excel.ActiveWorkbook

$setUp = $wb.Worksheets.Item("SetUp")
$loginPage = $wb.Worksheets.Item("LoginPage")
$migrationPage = $wb.Worksheets.Item("MigrationPage")

# New xpath rows on the LoginPage sheet.
$loginPage.Range("A8").Value = "PlansAndAddonsUsername"
$loginPage.Range("B8").Value = "tvsap199aa"
$loginPage.Range("A9").Value = "PlansAndAddonsPassword"
# Leading apostrophe forces text entry (quote-prefix), matching the other
# password rows (B3/B5/B7) which store "123456" as text, not a number.
$loginPage.Range("B9").Value = "'123456"

# Widen the key column now that it holds a longer string (target raw
# OOXML width is 24.88671875; the COM ColumnWidth setter only offers
# whole-pixel granularity, so 24 -> stored width 24.8333... is the
# closest reachable value).
$loginPage.Columns.Item(1).ColumnWidth = 24

# Move the selection on SetUp (not the active tab at the end, so select
# first) and then on LoginPage.
$setUp.Range("B18").Select()
$loginPage.Range("B27").Select()

# MigrationPage ends up the active tab/sheet.
$migrationPage.Activate()
